$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (substory_id) values for rows 3..13 to reflect the new data.
# Row 2's C value (3) is unchanged.
$ws.Range("C3").Value = 9
$ws.Range("C4").Value = 13
$ws.Range("C5").Value = 16
$ws.Range("C6").Value = 20
$ws.Range("C7").Value = 16
$ws.Range("C8").Value = 23
$ws.Range("C9").Value = 26
$ws.Range("C10").Value = 30
$ws.Range("C11").Value = 40
$ws.Range("C12").Value = 34
$ws.Range("C13").Value = 36

# Update the active selection to match the edited cell F6.
$ws.Range("F6").Select()
